$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new header cell BB1 (quarter date serial for 2025-11-25)
$ws.Cells.Item(1, 54).Value = 45986

# Rows 2-82: new column BB (column 54), mirroring/extending the BA column
$ws.Cells.Item(2, 54).Value = 1.63165117816655
$ws.Cells.Item(3, 54).Value = 1.358419700277324
$ws.Cells.Item(4, 54).Value = -0.1398451755242718
$ws.Cells.Item(5, 54).Value = -0.07003068004286206
$ws.Cells.Item(6, 54).Value = -7.3868590581191
$ws.Cells.Item(7, 54).Value = -10.33903131837766
$ws.Cells.Item(8, 54).Value = -1.060235412222937
$ws.Cells.Item(9, 54).Value = 3.423433284297019
$ws.Cells.Item(10, 54).Value = 2.063081011733999
$ws.Cells.Item(11, 54).Value = 2.99202665168275
$ws.Cells.Item(12, 54).Value = 6.534207423404695
$ws.Cells.Item(13, 54).Value = 1.294218807309846
$ws.Cells.Item(14, 54).Value = 2.007778863461724
$ws.Cells.Item(15, 54).Value = 2.488825190489734
$ws.Cells.Item(16, 54).Value = -0.02236378853282872
$ws.Cells.Item(17, 54).Value = 1.864148130528193
$ws.Cells.Item(18, 54).Value = 0.3074416423962276
$ws.Cells.Item(19, 54).Value = 1.591452302439862
$ws.Cells.Item(20, 54).Value = 1.24229439238384
$ws.Cells.Item(21, 54).Value = 1.021051258256691
$ws.Cells.Item(22, 54).Value = -2.083516682089652
$ws.Cells.Item(23, 54).Value = 0.8945075486844729
$ws.Cells.Item(24, 54).Value = 1.373553924371535
$ws.Cells.Item(25, 54).Value = 0.3276711086097635
$ws.Cells.Item(26, 54).Value = 1.809509533486136
$ws.Cells.Item(27, 54).Value = 0.4339813219297071
$ws.Cells.Item(28, 54).Value = 0.4041977135476031
$ws.Cells.Item(29, 54).Value = 1.507492882068462
$ws.Cells.Item(30, 54).Value = 1.350099422002103
$ws.Cells.Item(31, 54).Value = 1.114800352984872
$ws.Cells.Item(32, 54).Value = 1.516137977472326
$ws.Cells.Item(33, 54).Value = -0.02409266101658147
$ws.Cells.Item(34, 54).Value = -0.7388860519741201
$ws.Cells.Item(35, 54).Value = 1.567695497950282
$ws.Cells.Item(36, 54).Value = 1.091382109524247
$ws.Cells.Item(37, 54).Value = -0.2679292251141305
$ws.Cells.Item(38, 54).Value = 1.682996656587392
$ws.Cells.Item(39, 54).Value = 1.591260107191601
$ws.Cells.Item(40, 54).Value = 1.13706115148689
$ws.Cells.Item(41, 54).Value = 1.735453665039003
$ws.Cells.Item(42, 54).Value = 2.625533283765208
$ws.Cells.Item(43, 54).Value = -0.3243298885145123
$ws.Cells.Item(44, 54).Value = 0.6372341859553217
$ws.Cells.Item(45, 54).Value = -0.5789332341234967
$ws.Cells.Item(46, 54).Value = 0.1728728569232914
$ws.Cells.Item(47, 54).Value = 1.7
$ws.Cells.Item(48, 54).Value = -1.4
$ws.Cells.Item(49, 54).Value = 1.3
$ws.Cells.Item(50, 54).Value = -0.6
$ws.Cells.Item(51, 54).Value = -3.292009884772611
$ws.Cells.Item(52, 54).Value = -20.353166912592
$ws.Cells.Item(53, 54).Value = 17.87631681612835
$ws.Cells.Item(54, 54).Value = 4.395171409529297
$ws.Cells.Item(55, 54).Value = 1.763269038133103
$ws.Cells.Item(56, 54).Value = 1.651360288740889
$ws.Cells.Item(57, 54).Value = -0.6023499639201475
$ws.Cells.Item(58, 54).Value = 3.773751128807561
$ws.Cells.Item(59, 54).Value = 0.3486937127494798
$ws.Cells.Item(60, 54).Value = 0.6021692060546258
$ws.Cells.Item(61, 54).Value = 1.505090528591751
$ws.Cells.Item(62, 54).Value = -1.274647749701984
$ws.Cells.Item(63, 54).Value = 0.4647476577308112
$ws.Cells.Item(64, 54).Value = -0.6885297541760451
$ws.Cells.Item(65, 54).Value = -0.9272934801906558
$ws.Cells.Item(66, 54).Value = -0.8572418502401149
$ws.Cells.Item(67, 54).Value = 0.05364604092510206
$ws.Cells.Item(68, 54).Value = 1.830682919206694
$ws.Cells.Item(69, 54).Value = -2.604576060482884
$ws.Cells.Item(70, 54).Value = -3.124437332092583
$ws.Cells.Item(71, 54).Value = 2.411315004676197
$ws.Cells.Item(72, 54).Value = 0.2886831937783967
$ws.Cells.Item(73, 54).Value = -0.7196185376451893
$ws.Cells.Item(74, 54).Value = -0.4337369856241082
$ws.Cells.Item(75, 54).Value = -0.4337369856241082
$ws.Cells.Item(76, 54).Value = -0.4337369856241082
$ws.Cells.Item(77, 54).Value = -0.4337369856241082
$ws.Cells.Item(78, 54).Value = -0.4337369856241082
$ws.Cells.Item(79, 54).Value = -0.4337369856241082
$ws.Cells.Item(80, 54).Value = -0.4337369856241082
$ws.Cells.Item(81, 54).Value = -0.4337369856241082
$ws.Cells.Item(82, 54).Value = -0.4337369856241082

# New row 83: extra forecast quarter
# Copy the date-style formatting (style index 1) from A82 onto the new A83 cell
$ws.Cells.Item(82, 1).Copy() | Out-Null
$ws.Cells.Item(83, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(83, 1).Value = 46934
$ws.Cells.Item(83, 54).Value = -0.4337369856241082

$excel.CutCopyMode = 0

Write-Output "done"
